$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.029.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "'3.321.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.05%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'600.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "'144.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.38%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'3.319.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.26%  "
$ws.Range("D9").Value = "'0.524"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("D10").Value = "'0.150"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.08%  "
$ws.Range("D11").Value = "'5.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.51%  "
$ws.Range("D12").Value = "'0.475"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.84%  "
$ws.Range("D13").Value = "'0.0000251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("D14").Value = "'35.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.25%  "
$ws.Range("D15").Value = "'3.869.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.17%  "
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "'3.327.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.34%  "
$ws.Range("D18").Value = "'64.136.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").Value = "'6.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.23%  "
$ws.Range("D20").Value = "'484.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.34%  "
$ws.Range("D21").Value = "'14.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").Value = "'0.740"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.87%  "
$ws.Range("D23").Value = "'8.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.91%  "
$ws.Range("D24").Value = "'13.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.72%  "
$ws.Range("D25").Value = "'85.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.89%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").Value = "  +2.30%  "
$ws.Range("D28").Value = "'8.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.70%  "
$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'7.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.86%  "
$ws.Range("D31").Value = "'29.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.60%  "
$ws.Range("E32").Value = "  +5.47%  "
$ws.Range("D33").Value = "'0.106"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.04%  "
$ws.Range("D34").Value = "'2.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.63%  "
$ws.Range("E35").Value = "  +2.69%  "
$ws.Range("D36").Value = "'6.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.33%  "
$ws.Range("D37").Value = "'0.0₃0763"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.59%  "
$ws.Range("D38").Value = "'53.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("E39").Value = "  +4.61%  "
$ws.Range("D40").Value = "'437.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.29%  "
$ws.Range("D41").Value = "'3.057.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.67%  "
$ws.Range("D42").Value = "'2.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.83%  "
$ws.Range("D43").Value = "'8.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.56%  "
$ws.Range("E44").Value = "  -1.24%  "
$ws.Range("D45").Value = "'0.269"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.59%  "
$ws.Range("D46").Value = "'2.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.87%  "
$ws.Range("D47").Value = "'26.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.78%  "
$ws.Range("D48").Value = "'36.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +13.50%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("D51").Value = "'2.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.72%  "
